# Adds the new "Container With Most Water" Array problem as row 17,
# following the same Name/Description/Approach/Link layout as the
# existing rows (e.g. row 16 - "3Sum").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row content (A17:C17) -------------------------------------------------
$ws.Range("A17").Value2 = "Container With Most Water"
$ws.Range("B17").Value2 = "Return maximum possible area"
$ws.Range("C17").Value2 = "Use two pointers left and right to calculate current area. Compare height at left and right pointer. Move pointer inwards whichever has shorter height."

# --- Link cell (D17) -------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D17"), "https://leetcode.com/problems/container-with-most-water/")

# --- Match formatting of the row above (A16:D16 is the "Neutral" category row) --
$ws.Range("A16:D16").Copy()
$ws.Range("A17:D17").PasteSpecial(-4122)

# --- Restore the selection used by the author in the final state of the file ---
$ws.Range("C14").Select()
